$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR")
$ws.Rows.Item(10).Delete()
$excel.CalculateFullRebuild()
$ws.Activate()
